$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 21:57"

$ws.Range("B4").Value = 4746948
$ws.Range("C4").Value = 41059
$ws.Range("D4").Value = 2342867
$ws.Range("E4").Value = 2246688
$ws.Range("G4").Value = 646
$ws.Range("H4").Value = 157393

$ws.Range("B21").Value = 211060
$ws.Range("C21").Value = 395
$ws.Range("E21").Value = 8234
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 9226

$ws.Range("B31").Value = 86232
$ws.Range("C31").Value = 877
$ws.Range("D31").Value = 36213
$ws.Range("E31").Value = 44283
$ws.Range("G31").Value = 34
$ws.Range("H31").Value = 5736

$ws.Range("B65").Value = 24783
$ws.Range("C65").Value = 774
$ws.Range("D65").Value = 15299
$ws.Range("E65").Value = 9337
$ws.Range("G65").Value = 6
$ws.Range("H65").Value = 147

$ws.Range("A70").Value = "Costa Rica"
$ws.Range("B70").Value = 18187
$ws.Range("C70").Value = 367
$ws.Range("D70").Value = 4531
$ws.Range("E70").Value = 13502
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 154

$ws.Range("A71").Value = "Etiopia"
$ws.Range("B71").Value = 17999
$ws.Range("C71").Value = 469
$ws.Range("D71").Value = 7195
$ws.Range("E71").Value = 10520
$ws.Range("G71").Value = 10
$ws.Range("H71").Value = 284

$ws.Range("B76").Value = 16109
$ws.Range("C76").Value = 62
$ws.Range("D76").Value = 11750
$ws.Range("E76").Value = 4257

$ws.Range("B96").Value = 6319
$ws.Range("C96").Value = 9
$ws.Range("D96").Value = 5043
$ws.Range("E96").Value = 1119

$ws.Range("B104").Value = 4614
$ws.Range("C104").Value = 6
$ws.Range("D104").Value = 1635
$ws.Range("E104").Value = 2920

$ws.Range("B108").Value = 3949
$ws.Range("C108").Value = 156
$ws.Range("D108").Value = 2613
$ws.Range("E108").Value = 1320

$ws.Range("B119").Value = 2706
$ws.Range("C119").Value = 58
$ws.Range("E119").Value = 1449
$ws.Range("G119").Value = 2
$ws.Range("H119").Value = 43

$ws.Range("B122").Value = 2480
$ws.Range("C122").Value = 29
$ws.Range("D122").Value = 1837
$ws.Range("E122").Value = 619
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 24

$ws.Range("A123").Value = "Sudan del Sur"
$ws.Range("B123").Value = 2352
$ws.Range("C123").Value = 30
$ws.Range("D123").Value = 1175
$ws.Range("E123").Value = 1131
$ws.Range("H123").Value = 46

$ws.Range("A124").Value = "Eslovaquia"
$ws.Range("B124").Value = 2337
$ws.Range("C124").Value = 45
$ws.Range("D124").Value = 1742
$ws.Range("E124").Value = 566
$ws.Range("H124").Value = 29

$ws.Range("A125").Value = "Namibia"
$ws.Range("B125").Value = 2224
$ws.Range("C125").Value = 95
$ws.Range("D125").Value = 171
$ws.Range("E125").Value = 2042
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 11

$ws.Range("A126").Value = "Eslovenia"
$ws.Range("B126").Value = 2171
$ws.Range("C126").Value = 15
$ws.Range("D126").Value = 1821
$ws.Range("E126").Value = 231
$ws.Range("H126").Value = 119

$ws.Range("B145").Value = 1164
$ws.Range("C145").Value = 16
$ws.Range("D145").Value = 460
$ws.Range("E145").Value = 650
$ws.Range("G145").Value = 2
$ws.Range("H145").Value = 54

$ws.Range("A146").Value = "Burkina Faso"
$ws.Range("B146").Value = 1143
$ws.Range("C146").Value = 37
$ws.Range("D146").Value = 935
$ws.Range("E146").Value = 155
$ws.Range("H146").Value = 53

$ws.Range("A147").Value = "Niger"
$ws.Range("B147").Value = 1134
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 1028
$ws.Range("E147").Value = 37
$ws.Range("H147").Value = 69

$ws.Range("A148").Value = "Republica de Chipre"
$ws.Range("B148").Value = 1124
$ws.Range("C148").Value = 5
$ws.Range("D148").Value = 852
$ws.Range("E148").Value = 253
$ws.Range("H148").Value = 19

$ws.Range("B153").Value = 874
$ws.Range("C153").Value = 3
$ws.Range("D153").Value = 782
$ws.Range("E153").Value = 77

$ws.Range("B168").Value = 386
$ws.Range("C168").Value = 8
$ws.Range("E168").Value = 49
